$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B76").Copy()
$ws.Range("B78").PasteSpecial(-4122)
$ws.Range("B78").Value = @'
The NeHA Operator is responsible for managing the system's day-to-day operations, with a focus on pharmacists, pharmacies, and certification processes. The NeHA Operator has the following responsibilities: Viewing all pharmacists, physicians, and pharmacies within the system. The ability to edit details of pharmacists and pharmacies to ensure up-to-date and accurate information. Enabling or disabling pharmacists and physicians within the system, ensuring only authorized individuals have access. Viewing active, expired, and revoked certificates, along with the authority to revoke certificates, ensuring compliance and integrity in certification processes
'@
$ws.Range("C78").Value = "Correct"

$ws.Range("B76").Copy()
$ws.Range("B79").PasteSpecial(-4122)
$ws.Range("B79").Value = @'
The Pharmacy Operator plays a crucial role in the operational management of pharmacies and pharmacists within the system. They are responsible for adding new pharmacists and pharmacies to the system, managing the details of pharmacies, and enabling or disabling pharmacists. They also have the authority to create new certificates, view their active, revoked, and expired certificates, and revoke active certificates..
'@
$ws.Range("C79").Value = "Correct"

$ws.Range("B76").Copy()
$ws.Range("B80").PasteSpecial(-4122)
$ws.Range("B80").Value = @'
Certificates play a crucial role in the roles of pharmacists and physicians within the system. They are used to verify and authenticate professional qualifications and services, ensuring that only qualified and authorized individuals can access and use the system. Certificates are created and managed by pharmacists and physicians, who are responsible for ensuring that the information contained in the certificates is accurate and up-to-date. This includes information such as the individual's name, qualifications, and any restrictions or limitations on their practice. Certificates are also used to verify the authenticity of professional qualifications and services. For example, a pharmacist may use a certificate to verify the qualifications of a physician before prescribing a medication. Similarly, a physician may use a certificate to verify the qualifications of a pharmacist before dispensing a prescription. In summary, certificates are an essential tool for ensuring the quality and safety of healthcare services in the system. They provide a secure and reliable way to verify and authenticate professional qualifications and services, ensuring that only qualified and authorized individuals can access and use the system..
'@
$ws.Range("C80").Value = "Correct"

$ws.Range("B76").Copy()
$ws.Range("B81").PasteSpecial(-4122)
$ws.Range("B81").Value = @'
Physicians and pharmacists play a crucial role in the NCP eHealth(Cyprus) Portal regarding certificates. Physicians are responsible for creating and managing certificates, while pharmacists are authorized to view their active, revoked, and expired certificates. Both roles have the authority to revoke active certificates, ensuring the system reflects their current professional standing accurately..
'@
$ws.Range("C81").Value = "Correct"

$ws.Range("B76").Copy()
$ws.Range("B82").PasteSpecial(-4122)
$ws.Range("B82").Value = @'
The responsibilities associated with the creation of certificates for physicians and pharmacists are as follows: Physicians: - Create new certificates for themselves and other professionals. - View their active, revoked, and expired certificates. - Revocate active certificates. Pharmacists: - Create new certificates for themselves and other professionals.
'@
$ws.Range("C82").Value = "Correct"

$ws.Range("B76").Copy()
$ws.Range("B83").PasteSpecial(-4122)
$ws.Range("B83").Value = @'
Physicians and pharmacists can view their certificates within the system by logging into the NCPeH CY portal and navigating to the "My Certificates" section. This section provides a comprehensive overview of their active, revoked, and expired certificates, allowing them to monitor their certification status and ensure that their professional qualifications and services are accurately reflected in the system..
'@
$ws.Range("C83").Value = "Correct"

$ws.Range("B74").Copy()
$ws.Range("B84").PasteSpecial(-4122)
$ws.Range("B84").Value = @'
The NCP eHealth(Cyprus) Portal is a digital platform that provides a secure and efficient way for healthcare professionals to share and access patient information. The portal is designed to facilitate cross-border eHealth services, allowing healthcare providers to collaborate and coordinate care for patients who may be receiving treatment in multiple countries. The portal is managed by the National eHealth Authority (NeHA) in Cyprus, which is responsible for the development and implementation of the portal. The NeHA is a government agency that is tasked with promoting the use of eHealth technologies in the healthcare sector. The portal is designed to be user-friendly and accessible to healthcare professionals from all over the world. It provides a secure and encrypted environment for healthcare providers to share patient information, including medical records, test results, and other relevant data. The portal is also designed to be compliant with international standards and regulations, ensuring that patient data is protected and secure. This includes measures such as encryption, authentication, and access control, which help to prevent unauthorized access to patient information. The portal is also designed to be scalable, allowing it to accommodate the needs of healthcare providers from different countries and regions. This means that the portal can be adapted to meet the specific needs of different healthcare systems, ensuring that it remains relevant and useful for healthcare providers around the world. Overall, the NCP eHealth(Cyprus) Portal is a valuable tool for healthcare providers who are looking to improve the quality of care for their patients. It provides a secure and efficient way for healthcare professionals to share and access patient information, while also ensuring that patient data is protected and secure..
'@
$ws.Range("C84").Value = "Incorrect"

$ws.Range("B76").Copy()
$ws.Range("B85").PasteSpecial(-4122)
$ws.Range("B85").Value = @'
Certificates are crucial for verification and authentication within the portal because they provide a secure and reliable way to verify the identity and credentials of healthcare professionals. Certificates are issued by trusted authorities and contain information such as the professional's name, qualifications, and other relevant details. This information is used to verify the professional's identity and credentials, ensuring that they are qualified to provide healthcare services. Certificates also provide a way to authenticate the professional's identity and credentials, ensuring that the information contained in the certificate is accurate and up-to-date. This helps to ensure that healthcare professionals are providing safe and effective care to their patients..
'@
$ws.Range("C85").Value = "Correct"

$excel.CutCopyMode = $false

$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$ws.Range("C85").Select() | Out-Null
